$d = $word.ActiveDocument

$d.Content.Find.Execute("PRIMER TRABAJO: ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
